$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, shifting existing rows 19..168 down to 20..169
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new record's data
$ws.Range("A19").Value = 4
$ws.Range("B19").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C19").Value = "Los Lagos"
$ws.Range("D19").Value = 44537
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 100112017
$ws.Range("G19").Value = "Apio"
$ws.Range("H19").Value = "Americana (o)"
$ws.Range("I19").Value = "Segunda"
$ws.Range("J19").Value = 35
$ws.Range("K19").Value = 10000
$ws.Range("L19").Value = 10000
$ws.Range("M19").Value = 10000
$ws.Range("N19").Value = "`$/docena de matas"
$ws.Range("O19").Value = "Región de Coquimbo"
$ws.Range("P19").Value = 1667
$ws.Range("Q19").Value = 6
$ws.Range("R19").Value = "Hortaliza"
